$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Consolidated age-range data: the "0 to 9", "10 to 19" and "20 to 29"
# buckets are merged into a single "0 to 29" bucket (summed), per region.
$data = @(
    @('Arecibo', '0 to 29', 5),
    @('Arecibo', '30 to 39', 7),
    @('Arecibo', '40 to 49', 23),
    @('Arecibo', '50 to 59', 51),
    @('Arecibo', '60 to 69', 104),
    @('Arecibo', '70 to 79', 139),
    @('Arecibo', '80+', 155),
    @('Arecibo', 'N/A', 0),
    @('Bayamon', '0 to 29', 7),
    @('Bayamon', '30 to 39', 18),
    @('Bayamon', '40 to 49', 39),
    @('Bayamon', '50 to 59', 117),
    @('Bayamon', '60 to 69', 162),
    @('Bayamon', '70 to 79', 220),
    @('Bayamon', '80+', 295),
    @('Bayamon', 'N/A', 0),
    @('Caguas', '0 to 29', 3),
    @('Caguas', '30 to 39', 11),
    @('Caguas', '40 to 49', 29),
    @('Caguas', '50 to 59', 81),
    @('Caguas', '60 to 69', 97),
    @('Caguas', '70 to 79', 155),
    @('Caguas', '80+', 212),
    @('Caguas', 'N/A', 0),
    @('Fajardo', '0 to 29', 5),
    @('Fajardo', '30 to 39', 4),
    @('Fajardo', '40 to 49', 12),
    @('Fajardo', '50 to 59', 34),
    @('Fajardo', '60 to 69', 33),
    @('Fajardo', '70 to 79', 40),
    @('Fajardo', '80+', 52),
    @('Fajardo', 'N/A', 0),
    @('Mayaguez', '0 to 29', 3),
    @('Mayaguez', '30 to 39', 10),
    @('Mayaguez', '40 to 49', 26),
    @('Mayaguez', '50 to 59', 51),
    @('Mayaguez', '60 to 69', 88),
    @('Mayaguez', '70 to 79', 145),
    @('Mayaguez', '80+', 187),
    @('Mayaguez', 'N/A', 0),
    @('Metro', '0 to 29', 13),
    @('Metro', '30 to 39', 27),
    @('Metro', '40 to 49', 79),
    @('Metro', '50 to 59', 143),
    @('Metro', '60 to 69', 203),
    @('Metro', '70 to 79', 282),
    @('Metro', '80+', 388),
    @('Metro', 'N/A', 0),
    @('Ponce', '0 to 29', 4),
    @('Ponce', '30 to 39', 6),
    @('Ponce', '40 to 49', 14),
    @('Ponce', '50 to 59', 27),
    @('Ponce', '60 to 69', 52),
    @('Ponce', '70 to 79', 86),
    @('Ponce', '80+', 103),
    @('Ponce', 'N/A', 0),
    @('N/A', '0 to 29', 1),
    @('N/A', '30 to 39', 1),
    @('N/A', '40 to 49', 4),
    @('N/A', '50 to 59', 15),
    @('N/A', '60 to 69', 14),
    @('N/A', '70 to 79', 11),
    @('N/A', '80+', 6)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $data[$i]
    $r = $i + 2
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
}

# The old sheet had 80 data rows (81 incl. header); the consolidated sheet
# only needs 63 data rows (through row 64). Remove the
# now-stale trailing rows.
$ws.Range("A65:C81").EntireRow.Delete()

# Restore view state (scroll position / zoom / selection) from the diff.
$ws.Application.ActiveWindow.Zoom = 161
$ws.Range("A45").Select()
$ws.Range("C67").Select()
